# Generate Report for Handoff
# Regenerates the localization status report: the source markdown file was
# renamed to a new GUID, and the latest handoff/handback timestamps advance.

$wb = $excel.ActiveWorkbook

# Old / new identifiers
$oldGuid = "5ed9a9ff-9c92-478c-873a-7fe0e3bccb6e"
$newGuid = "a5cfe79f-3c5f-4df4-8006-b4782083b8ec"

$oldHash = "83b46e55fa0e936070718a24612133eea82259e9"
$newHash = "a0fe7ee353b6f86ce7cb473353175e76f8d7bec5"

$newHoGenerateDate    = "2016-09-07 09:29:50"
$newZhHandbackDate    = "2016-09-07 09:29:44"

$newMdName  = "$newGuid.md"
$newMdDisp  = "e2e\$newGuid.md"

# The hyperlink's underlying target URL is untouched by this edit - only the
# visible display text changes - so keep the existing Address for each link.
$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7607fc90cc44f1058a85538563a416dde979a4a7/e2e/$oldGuid.md"

function Update-HyperlinkDisplay($range, $url, $displayText) {
    # TextToDisplay can't be edited in place on this host (it always appends
    # a new hyperlink instead of updating the existing one), so drop the
    # existing hyperlink on the cell and re-add it with the same address but
    # the new display text.
    $range.Hyperlinks.Delete()
    $range.Worksheet.Hyperlinks.Add($range, $url, $null, $null, $displayText) | Out-Null
}

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
Update-HyperlinkDisplay $wsOverview.Range("B2") $hyperlinkUrl $newMdDisp
$wsOverview.Range("G2").Value = $newHoGenerateDate

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
Update-HyperlinkDisplay $wsZh.Range("A2") $hyperlinkUrl $newMdName
$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = $newZhHandbackDate

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
Update-HyperlinkDisplay $wsDe.Range("A2") $hyperlinkUrl $newMdName
$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = $newHoGenerateDate
